$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (id = 1)
$ws.Range("B2").Value = 3.601030655423489
$ws.Range("D2").Value = 6.373286323410602
$ws.Range("F2").Value = 3.983303952131626
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 2.717403253259548
$ws.Range("L2").Value = 1.330129092892181

# Row 3 (id = 2)
$ws.Range("B3").Value = 2.550588267952786
$ws.Range("D3").Value = 2.711194164626598
$ws.Range("H3").Value = 2.389982371278975

# Row 4 (id = 3)
$ws.Range("B4").Value = 4.282329885966153
$ws.Range("D4").Value = 6.174801091503967
$ws.Range("F4").Value = 2.389858680428339
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
